$d = $word.ActiveDocument

# Original text is a single run "hello". Target splits it into three runs:
#   "H"  +  "ello"  +  " this is first line"
# all sharing the same run formatting (lang=en-AU), so the split has to be
# forced explicitly (the engine otherwise re-coalesces adjacent runs that
# end up with identical formatting).

# Step 1: capitalize the leading "h" -> "H" (pure text edit, no split yet).
$rH = $d.Range(0, 1)
$rH.Text = "H"

# Step 2: insert the new sentence right after "hello" / before the
# _GoBack bookmark. Assigning .Text on a collapsed range (rather than
# InsertAfter) keeps the inserted text ahead of the bookmark markers.
$rTail = $d.Range(5, 5)
$rTail.Text = " this is first line"

# Step 3: force the run boundaries the diff expects by toggling a
# formatting property on and off for each sub-range; toggling back to the
# original value keeps the visible formatting unchanged while still
# preventing the engine from merging the runs back together.
$rH2 = $d.Range(0, 1)
$rH2.Bold = 1
$rH2.Bold = 0

$rNew = $d.Range(5, 25)
$rNew.Bold = 1
$rNew.Bold = 0

Write-Output $d.Content.Text
